$d = $word.ActiveDocument

$replacements = @(
    @("53×18=", "19×27="),
    @("45×25=", "51×87="),
    @("76×18=", "97×34="),
    @("72×77=", "62×82="),
    @("43×72=", "75×87="),
    @("75×84=", "38×35="),
    @("38×61=", "82×67="),
    @("70×47=", "85×67="),
    @("56×61=", "62×60="),
    @("78×74=", "66×87="),
    @("61×83=", "84×30="),
    @("72×63=", "60×56="),
    @("72×22=", "95×82="),
    @("31×46=", "63×71="),
    @("54×46=", "71×78="),
    @("35×20=", "81×28="),
    @("38×56=", "65×30="),
    @("91×21=", "44×34="),
    @("72×17=", "60×41="),
    @("99×70=", "61×71="),
    @("49×16=", "66×81="),
    @("42×42=", "76×41="),
    @("28×89=", "44×93="),
    @("57×91=", "23×23="),
    @("87×18=", "73×45=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
